# Apply metadata updates to the "Metadata" worksheet to reflect the
# repository move from John Moehrke's personal GitHub Pages site to the
# Department of Veterans Affairs GitHub org, per the commit message:
# "changing to be more formal now that it is VA github repo"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL: johnmoehrke.github.io/MHV-PHR -> department-of-veterans-affairs.github.io/mhv-fhir-phr-mapping
$ws.Range("B2").Value = "https://department-of-veterans-affairs.github.io/mhv-fhir-phr-mapping/ValueSet/AllergyCategoryVS"

# Version: 0.1.20-beta -> 0.2.0
$ws.Range("B3").Value = "0.2.0"

# Date: 2023-08-22T12:58:52-05:00 -> 2023-08-22T16:36:15-05:00
$ws.Range("B8").Value = "2023-08-22T16:36:15-05:00"

# Publisher: John Moehrke (himself) -> VA Digital Services
$ws.Range("B9").Value = "VA Digital Services"
